$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# The table currently has 35 columns: id_region, id_sector, id_subsector, unit, 2020..2050.
# Add 10 new columns (they are appended at the end of the table by this host).
for ($i = 0; $i -lt 10; $i++) {
    $null = $tbl.ListColumns.Add()
}

# Relabel the header row so that the year columns read 2010..2050 in order,
# right after the fixed id_region/id_sector/id_subsector/unit columns.
$years = 2010..2050
$headerRange = $tbl.HeaderRowRange
for ($i = 0; $i -lt $years.Count; $i++) {
    $headerRange.Cells.Item(1, 4 + $i + 1).Value = [string]$years[$i]
}

# All year columns carry the value 1 for every data row (same as the pre-existing years).
$dataRange = $ws.Range("E2:AS18")
$dataRange.Value = 1

$ws.Range("H10").Select()
